# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.880.13'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.744.48'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''573.06'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '''156.57'
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").Value = '''0.599'
$ws.Range("E8").Value = '  -1.53%  '
$ws.Range("D9").Value = '''0.109'
$ws.Range("E9").Value = '  -2.88%  '
$ws.Range("E10").Value = '  +0.22%  '
$ws.Range("D11").Value = '''0.381'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("D12").Value = '''5.33'
$ws.Range("E12").Value = '  -20.70%  '
$ws.Range("D13").Value = '3.227.21'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").Value = '''26.42'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").Value = '63.560.61'
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("D16").Value = '''0.0000149'
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").Value = '2.747.11'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '''12.08'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("D19").Value = '''4.78'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").Value = '''353.64'
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").Value = '''6.72'
$ws.Range("E21").Value = '  -3.34%  '
$ws.Range("D22").Value = '''0.997'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''0.534'
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("D24").Value = '''64.74'
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("E25").Value = '  -1.75%  '
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = '''8.32'
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("D28").Value = '0.0₃0885'
$ws.Range("E28").Value = '  -2.84%  '
$ws.Range("D29").Value = '''1.92'
$ws.Range("E29").Value = '  -3.98%  '
$ws.Range("D30").Value = '''6.88'
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("D31").Value = '''169.41'
$ws.Range("E31").Value = '  -1.48%  '
$ws.Range("D32").Value = '''1.21'
$ws.Range("E32").Value = '  -4.63%  '
$ws.Range("D33").Value = '''20.03'
$ws.Range("E33").Value = '  -1.97%  '
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").Value = '''1.77'
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("D38").Value = '''0.971'
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("D39").Value = '''6.13'
$ws.Range("E39").Value = '  +7.66%  '
$ws.Range("D40").Value = '''4.10'
$ws.Range("E40").Value = '  -2.93%  '
$ws.Range("D41").Value = '''323.17'
$ws.Range("E41").Value = '  -6.24%  '
$ws.Range("D42").Value = '''38.90'
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '''21.08'
$ws.Range("E43").Value = '  -3.45%  '
$ws.Range("D44").Value = '''0.0583'
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = '''21.14'
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("D46").Value = '''134.77'
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0252'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.623'
$ws.Range("E48").Value = '  -3.86%  '
$ws.Range("D49").Value = '''0.100'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").Value = '''0.999'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("D51").Value = '''11.04'
$ws.Range("E51").Value = '  +0.50%  '
